$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140..259 down to 141..260
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new weekly price record
$ws.Range("A140").Value = 5
$ws.Range("B140").Value = "Macroferia Regional de Talca"
$ws.Range("C140").Value = "Maule"
$ws.Range("D140").Value = 45096
$ws.Range("E140").Value = 7
$ws.Range("F140").Value = 100112031
$ws.Range("G140").Value = "Poroto verde"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 150
$ws.Range("K140").Value = 23000
$ws.Range("L140").Value = 23000
$ws.Range("M140").Value = 23000
$ws.Range("N140").Value = "`$/malla 25 kilos"
$ws.Range("O140").Value = "Perú"
$ws.Range("P140").Value = 920
$ws.Range("Q140").Value = 25
$ws.Range("R140").Value = "Hortaliza"
